$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6945960
$ws.Range("I19").Value = 13889399
$ws.Range("J19").Value = 2521.3333
$ws.Range("K19").Value = 13889399
$ws.Range("L19").Value = 2521.3333
$ws.Range("M19").Value = -13889224
$ws.Range("N19").Value = -2871.3333
$ws.Range("H28").Value = 25002358
$ws.Range("I28").Value = 31252804
$ws.Range("J28").Value = 575
$ws.Range("K28").Value = 31252804
$ws.Range("L28").Value = 575
$ws.Range("M28").Value = -31252319
$ws.Range("N28").Value = -1545
$ws.Range("H32").Value = 752.5833
$ws.Range("I32").Value = 697.625
$ws.Range("J32").Value = 862.5
$ws.Range("K32").Value = 697.625
$ws.Range("L32").Value = 862.5
$ws.Range("M32").Value = -371.625
$ws.Range("N32").Value = -1514.5
$ws.Range("H33").Value = 6363
$ws.Range("I33").Value = 5150
$ws.Range("J33").Value = 10002
$ws.Range("K33").Value = 5150
$ws.Range("L33").Value = 10002
$ws.Range("M33").Value = -4921
$ws.Range("N33").Value = -10460
$ws.Range("H53").Value = 9655.9375
$ws.Range("I53").Value = 25140
$ws.Range("J53").Value = 365.5
$ws.Range("K53").Value = 25140
$ws.Range("L53").Value = 365.5
$ws.Range("M53").Value = -24503
$ws.Range("N53").Value = -1639.5
$ws.Range("H62").Value = 116670580
$ws.Range("I62").Value = 50005400
$ws.Range("J62").Value = 250000940
$ws.Range("K62").Value = 50005400
$ws.Range("L62").Value = 250000940
$ws.Range("M62").Value = -50004776
$ws.Range("N62").Value = -250002188
$ws.Range("H65").Value = 116670580
$ws.Range("I65").Value = 50005400
$ws.Range("J65").Value = 250000940
$ws.Range("K65").Value = 250027000
$ws.Range("L65").Value = 1250004700
$ws.Range("M65").Value = -250023880
$ws.Range("N65").Value = -1250010940
$ws.Range("H100").Value = 15386322
$ws.Range("I100").Value = 875.4286
$ws.Range("J100").Value = 33336008
$ws.Range("K100").Value = 875.4286
$ws.Range("L100").Value = 33336008
$ws.Range("M100").Value = -334.4286
$ws.Range("N100").Value = -33337090
$ws.Range("H107").Value = 1134.3077
$ws.Range("I107").Value = 1239.35
$ws.Range("J107").Value = 784.1667
$ws.Range("K107").Value = 1239.35
$ws.Range("L107").Value = 784.1667
$ws.Range("M107").Value = 680.6500000000001
$ws.Range("N107").Value = -4624.1667
$ws.Range("H113").Value = 4169008.5
$ws.Range("I113").Value = 9093019
$ws.Range("J113").Value = 2538.4614
$ws.Range("K113").Value = 9093019
$ws.Range("L113").Value = 2538.4614
$ws.Range("M113").Value = -9089765
$ws.Range("N113").Value = -9046.4614
$ws.Range("H132").Value = 2317421.8
$ws.Range("I132").Value = 3118.457
$ws.Range("J132").Value = 8548239
$ws.Range("K132").Value = 9355.370999999999
$ws.Range("L132").Value = 25644717
$ws.Range("M132").Value = -6825.370999999999
$ws.Range("N132").Value = -25649777
$ws.Range("H141").Value = 3100
$ws.Range("I141").Value = 2486.875
$ws.Range("J141").Value = 5552.5
$ws.Range("K141").Value = 7460.625
$ws.Range("L141").Value = 16657.5
$ws.Range("M141").Value = -2280.625
$ws.Range("N141").Value = -27017.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1107.8948
$ws.Range("I122").Value = 975.6667
$ws.Range("J122").Value = 1226.9
$ws.Range("K122").Value = 2927.0001
$ws.Range("L122").Value = 3680.7
$ws.Range("M122").Value = -477.0001000000002
$ws.Range("N122").Value = -8580.700000000001
$ws.Range("H135").Value = 49355.266
$ws.Range("J135").Value = 49355.266
$ws.Range("L135").Value = 49355.266
$ws.Range("N135").Value = -59495.266
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 38482296
$ws.Range("I20").Value = 55567824
$ws.Range("J20").Value = 39862.25
$ws.Range("K20").Value = 55567824
$ws.Range("L20").Value = 39862.25
$ws.Range("M20").Value = -55567577
$ws.Range("N20").Value = -40356.25
$ws.Range("H80").Value = 259.58823
$ws.Range("I80").Value = 169.4
$ws.Range("J80").Value = 297.16666
$ws.Range("K80").Value = 169.4
$ws.Range("L80").Value = 297.16666
$ws.Range("M80").Value = 828.6
$ws.Range("N80").Value = -2293.16666
$ws.Range("H83").Value = 259.58823
$ws.Range("I83").Value = 169.4
$ws.Range("J83").Value = 297.16666
$ws.Range("K83").Value = 847
$ws.Range("L83").Value = 1485.8333
$ws.Range("M83").Value = 4145
$ws.Range("N83").Value = -11469.8333
$ws.Range("H134").Value = 10118335
$ws.Range("I134").Value = 11364590
$ws.Range("J134").Value = 4025533
$ws.Range("K134").Value = 34093770
$ws.Range("L134").Value = 12076599
$ws.Range("M134").Value = -34091235
$ws.Range("N134").Value = -12081669
$ws.Range("H135").Value = 45330
$ws.Range("J135").Value = 45330
$ws.Range("L135").Value = 45330
$ws.Range("N135").Value = -55470
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5000225.5
$ws.Range("I22").Value = 6250188
$ws.Range("J22").Value = 375
$ws.Range("K22").Value = 6250188
$ws.Range("L22").Value = 375
$ws.Range("M22").Value = -6249838
$ws.Range("N22").Value = -1075
$ws.Range("H58").Value = 1899102.8
$ws.Range("I58").Value = 8023.0713
$ws.Range("J58").Value = 4546614.5
$ws.Range("K58").Value = 8023.0713
$ws.Range("L58").Value = 4546614.5
$ws.Range("M58").Value = -7820.0713
$ws.Range("N58").Value = -4547020.5
$ws.Range("H94").Value = 29417802
$ws.Range("I94").Value = 1318.2858
$ws.Range("J94").Value = 50009340
$ws.Range("K94").Value = 1318.2858
$ws.Range("L94").Value = 50009340
$ws.Range("M94").Value = -867.2858000000001
$ws.Range("N94").Value = -50010242
$ws.Range("H99").Value = 20622.727
$ws.Range("I99").Value = 26666.666
$ws.Range("J99").Value = 18356.25
$ws.Range("K99").Value = 26666.666
$ws.Range("L99").Value = 18356.25
$ws.Range("M99").Value = -25168.666
$ws.Range("N99").Value = -21352.25
$ws.Range("H126").Value = 20622.727
$ws.Range("I126").Value = 26666.666
$ws.Range("J126").Value = 18356.25
$ws.Range("K126").Value = 79999.99800000001
$ws.Range("L126").Value = 55068.75
$ws.Range("M126").Value = -77529.99800000001
$ws.Range("N126").Value = -60008.75
$ws.Range("H132").Value = 2438.0454
$ws.Range("I132").Value = 1936.1111
$ws.Range("J132").Value = 2785.5386
$ws.Range("K132").Value = 5808.3333
$ws.Range("L132").Value = 8356.6158
$ws.Range("M132").Value = -3278.3333
$ws.Range("N132").Value = -13416.6158
$ws.Range("H134").Value = 1545311
$ws.Range("I134").Value = 7904.1763
$ws.Range("K134").Value = 23712.5289
$ws.Range("M134").Value = -21177.5289
$ws.Range("H136").Value = 1899102.8
$ws.Range("I136").Value = 8023.0713
$ws.Range("J136").Value = 4546614.5
$ws.Range("K136").Value = 24069.2139
$ws.Range("L136").Value = 13639843.5
$ws.Range("M136").Value = -21519.2139
$ws.Range("N136").Value = -13644943.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 66668980
$ws.Range("J69").Value = 90911440
$ws.Range("L69").Value = 272734320
$ws.Range("N69").Value = -272735942
$ws.Range("H72").Value = 66668980
$ws.Range("J72").Value = 90911440
$ws.Range("L72").Value = 818202960
$ws.Range("N72").Value = -818211072
$ws.Range("H132").Value = 2755.0908
$ws.Range("I132").Value = 3004
$ws.Range("J132").Value = 2730.2
$ws.Range("K132").Value = 27036
$ws.Range("L132").Value = 24571.8
$ws.Range("M132").Value = -24506
$ws.Range("N132").Value = -29631.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5742.6665
$ws.Range("I126").Value = 7791.647
$ws.Range("J126").Value = 2259.4
$ws.Range("K126").Value = 23374.941
$ws.Range("L126").Value = 6778.200000000001
$ws.Range("M126").Value = -20904.941
$ws.Range("N126").Value = -11718.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1292.5
$ws.Range("I7").Value = 1325
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 1325
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -1213
$ws.Range("N7").Value = -1224
$ws.Range("H46").Value = 872.2143
$ws.Range("I46").Value = 803.9722
$ws.Range("J46").Value = 1281.6666
$ws.Range("K46").Value = 803.9722
$ws.Range("L46").Value = 1281.6666
$ws.Range("M46").Value = -615.9722
$ws.Range("N46").Value = -1657.6666
$ws.Range("H61").Value = 2708.25
$ws.Range("I61").Value = 1433.1666
$ws.Range("J61").Value = 3983.3333
$ws.Range("K61").Value = 1433.1666
$ws.Range("L61").Value = 3983.3333
$ws.Range("M61").Value = -1231.1666
$ws.Range("N61").Value = -4387.3333
$ws.Range("H113").Value = 2708.25
$ws.Range("I113").Value = 1433.1666
$ws.Range("J113").Value = 3983.3333
$ws.Range("K113").Value = 1433.1666
$ws.Range("L113").Value = 3983.3333
$ws.Range("M113").Value = 736.8334
$ws.Range("N113").Value = -8323.3333
$ws.Range("H122").Value = 8197981
$ws.Range("I122").Value = 969658.9
$ws.Range("J122").Value = 40002600
$ws.Range("K122").Value = 2908976.7
$ws.Range("L122").Value = 120007800
$ws.Range("M122").Value = -2906526.7
$ws.Range("N122").Value = -120012700
$ws.Range("H126").Value = 1292.5
$ws.Range("I126").Value = 1325
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 3975
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -1505
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 4203665
$ws.Range("I132").Value = 6212291.5
$ws.Range("J132").Value = 3808.9092
$ws.Range("K132").Value = 18636874.5
$ws.Range("L132").Value = 11426.7276
$ws.Range("M132").Value = -18634344.5
$ws.Range("N132").Value = -16486.7276
$ws.Range("H136").Value = 5052613.5
$ws.Range("I136").Value = 5557814.5
$ws.Range("J136").Value = 599.5
$ws.Range("K136").Value = 16673443.5
$ws.Range("L136").Value = 1798.5
$ws.Range("M136").Value = -16670893.5
$ws.Range("N136").Value = -6898.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 284.83334
$ws.Range("I113").Value = 246.95238
$ws.Range("J113").Value = 550
$ws.Range("K113").Value = 740.8571400000001
$ws.Range("L113").Value = 1650
$ws.Range("M113").Value = 1429.14286
$ws.Range("N113").Value = -5990
$ws.Range("H132").Value = 1363226.2
$ws.Range("I132").Value = 2439.2964
$ws.Range("J132").Value = 5955882
$ws.Range("K132").Value = 7317.889200000001
$ws.Range("L132").Value = 17867646
$ws.Range("M132").Value = -4787.889200000001
$ws.Range("N132").Value = -17872706
$ws.Range("H136").Value = 5244.4517
$ws.Range("I136").Value = 1729.7778
$ws.Range("J136").Value = 6682.273
$ws.Range("K136").Value = 5189.3334
$ws.Range("L136").Value = 20046.819
$ws.Range("M136").Value = -1729.7778
$ws.Range("N136").Value = -25146.819
